$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures ---
$ws.Range("E11").Value = 88527      # VALOR MORA
$ws.Range("C13").Value = 1          # Cant. Trabajadores
$ws.Range("F13").Value = 3          # Cant. Periodos

# --- Row 16: replace MARIA BERNARDA SILVA FUENTES / period 2507 data with
#     LUIS ALBERTO PUELLO CASTELLON / period 1710 data ---
$ws.Range("C16").Value = "1047374641"
$ws.Range("D16").Value = "LUIS ALBERTO PUELLO CASTELLON"
$ws.Range("E16").Value = "1710"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 877803

# --- Row 17: same worker, period becomes 1711 ---
$ws.Range("E17").Value = "1711"

# --- Row 18: same worker, period becomes 1712 ---
$ws.Range("E18").Value = "1712"

# Row 18 becomes the last data row, so it should take on row 19's
# "closing border" formatting before row 19 is removed.
$ws.Range("B19:J19").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# --- Remove the now-duplicate last row (shifts the signature block up) ---
$ws.Rows("19").Delete()
